# Fix the "Name" header: it was previously stored with a trailing space
# ("Name "); re-enter it without the trailing space.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"

# Move the active selection, matching the saved cursor position.
$ws.Range("C5").Select()
